$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Head Plate (metal)"

$ws.Range("A31").Value = "Rubber Seal"
$ws.Range("B31").Value = 1
$ws.Hyperlinks.Add($ws.Range("C31"), "https://www.grainger.com/product/Buna-N-Round-Cord-Buna-N-6RTT7")
$ws.Range("C31").HorizontalAlignment = -4131

$ws.Range("H15").Select()
